# Daily attendance processing - 2025-10-31 04:26:03
# Reorders the "Recorded By" (column G) author lists so that any
# System-originated entries ("System" / "system") are moved to the end
# of the comma-separated list, preserving the relative order of the
# remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }
    if ($value -eq "") { continue }

    $parts = $value -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $nonSystem = @()
    $systemOnly = @()
    foreach ($p in $trimmed) {
        if ($p.ToLower() -eq "system") {
            $systemOnly += $p
        } else {
            $nonSystem += $p
        }
    }

    $newParts = $nonSystem + $systemOnly
    $newValue = [string]::Join(", ", $newParts)

    if ($newValue -ne $value) {
        $cell.Value2 = $newValue
    }
}
